$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$cs = $s.ColorScheme
$cs.Colors(1).RGB = 123456
Write-Output ("done")
